{"js": "// Replace the date line and the 26 division problems in the table with\n// their new values, as described by the diff.\nconst replacements = [\n  [\"2025-11-05 Wednesday\", \"2025-11-06 Thursday\"],\n  [\"379\u00f76=\", \"245\u00f79=\"],\n  [\"888\u00f72=\", \"482\u00f77=\"],\n  [\"665\u00f75=\", \"482\u00f75=\"],\n  [\"474\u00f74=\", \"875\u00f76=\"],\n  [\"117\u00f79=\", \"909\u00f72=\"],\n  [\"805\u00f77=\", \"221\u00f74=\"],\n  [\"854\u00f76=\", \"380\u00f72=\"],\n  [\"235\u00f79=\", \"264\u00f74=\"],\n  [\"511\u00f73=\", \"544\u00f78=\"],\n  [\"267\u00f75=\", \"575\u00f74=\"],\n  [\"525\u00f74=\", \"944\u00f79=\"],\n  [\"347\u00f79=\", \"549\u00f75=\"],\n  [\"245\u00f72=\", \"168\u00f75=\"],\n  [\"422\u00f76=\", \"480\u00f76=\"],\n  [\"314\u00f74=\", \"413\u00f73=\"],\n  [\"556\u00f76=\", \"889\u00f78=\"],\n  [\"324\u00f72=\", \"761\u00f76=\"],\n  [\"222\u00f79=\", \"677\u00f76=\"],\n  [\"431\u00f78=\", \"451\u00f79=\"],\n  [\"462\u00f77=\", \"691\u00f77=\"],\n  [\"401\u00f76=\", \"842\u00f72=\"],\n  [\"355\u00f77=\", \"411\u00f73=\"],\n  [\"839\u00f76=\", \"288\u00f79=\"],\n  [\"503\u00f74=\", \"319\u00f72=\"],\n  [\"659\u00f72=\", \"452\u00f77=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 26 division problems in the table with\n# their new values, as described by the diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-11-05 Wednesday\", \"2025-11-06 Thursday\"),\n    @(\"379\u00f76=\", \"245\u00f79=\"),\n    @(\"888\u00f72=\", \"482\u00f77=\"),\n    @(\"665\u00f75=\", \"482\u00f75=\"),\n    @(\"474\u00f74=\", \"875\u00f76=\"),\n    @(\"117\u00f79=\", \"909\u00f72=\"),\n    @(\"805\u00f77=\", \"221\u00f74=\"),\n    @(\"854\u00f76=\", \"380\u00f72=\"),\n    @(\"235\u00f79=\", \"264\u00f74=\"),\n    @(\"511\u00f73=\", \"544\u00f78=\"),\n    @(\"267\u00f75=\", \"575\u00f74=\"),\n    @(\"525\u00f74=\", \"944\u00f79=\"),\n    @(\"347\u00f79=\", \"549\u00f75=\"),\n    @(\"245\u00f72=\", \"168\u00f75=\"),\n    @(\"422\u00f76=\", \"480\u00f76=\"),\n    @(\"314\u00f74=\", \"413\u00f73=\"),\n    @(\"556\u00f76=\", \"889\u00f78=\"),\n    @(\"324\u00f72=\", \"761\u00f76=\"),\n    @(\"222\u00f79=\", \"677\u00f76=\"),\n    @(\"431\u00f78=\", \"451\u00f79=\"),\n    @(\"462\u00f77=\", \"691\u00f77=\"),\n    @(\"401\u00f76=\", \"842\u00f72=\"),\n    @(\"355\u00f77=\", \"411\u00f73=\"),\n    @(\"839\u00f76=\", \"288\u00f79=\"),\n    @(\"503\u00f74=\", \"319\u00f72=\"),\n    @(\"659\u00f72=\", \"452\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n}\n"}
